$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 15
$ws.Range("H15").Value = 1468.3125
$ws.Range("I15").Value = 1468.3125
$ws.Range("K15").Value = 4404.9375
$ws.Range("M15").Value = -4235.9375
# row 96
$ws.Range("H96").Value = 1107.8334
$ws.Range("I96").Value = 1231.125
$ws.Range("J96").Value = 861.25
$ws.Range("K96").Value = 3693.375
$ws.Range("L96").Value = 2583.75
$ws.Range("M96").Value = -2320.375
$ws.Range("N96").Value = -5329.75
# row 100
$ws.Range("H100").Value = 2754.8
$ws.Range("I100").Value = 1507.2858
$ws.Range("J100").Value = 5665.6665
$ws.Range("K100").Value = 1507.2858
$ws.Range("L100").Value = 5665.6665
$ws.Range("M100").Value = -966.2858000000001
$ws.Range("N100").Value = -6747.6665
# row 103
$ws.Range("H103").Value = 1949.3077
$ws.Range("I103").Value = 458
$ws.Range("J103").Value = 3227.5715
$ws.Range("K103").Value = 1374
$ws.Range("L103").Value = 9682.7145
$ws.Range("M103").Value = -788
$ws.Range("N103").Value = -10854.7145
# row 132
$ws.Range("H132").Value = 1275.9445
$ws.Range("I132").Value = 1091.6875
$ws.Range("K132").Value = 3275.0625
$ws.Range("M132").Value = -745.0625
# row 139
$ws.Range("H139").Value = 68570.57000000001
$ws.Range("J139").Value = 68570.57000000001
$ws.Range("L139").Value = 68570.57000000001
$ws.Range("N139").Value = -78850.57000000001

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 26
$ws.Range("H26").Value = 3169
$ws.Range("I26").Value = 3169
$ws.Range("K26").Value = 3169
$ws.Range("M26").Value = -2839
# row 32
$ws.Range("H32").Value = 3645.305
$ws.Range("I32").Value = 2870.7693
$ws.Range("K32").Value = 2870.7693
$ws.Range("M32").Value = -2583.7693
# row 61
$ws.Range("H61").Value = 6784.091
$ws.Range("I61").Value = 6292.2
$ws.Range("J61").Value = 7194
$ws.Range("K61").Value = 6292.2
$ws.Range("L61").Value = 7194
$ws.Range("M61").Value = -6080.2
$ws.Range("N61").Value = -7618
# row 74
$ws.Range("H74").Value = 17547268
$ws.Range("I74").Value = 37040280
$ws.Range("J74").Value = 3554
$ws.Range("K74").Value = 37040280
$ws.Range("L74").Value = 3554
$ws.Range("M74").Value = -37039406
$ws.Range("N74").Value = -5302
# row 77
$ws.Range("H77").Value = 17547268
$ws.Range("I77").Value = 37040280
$ws.Range("J77").Value = 3554
$ws.Range("K77").Value = 185201400
$ws.Range("L77").Value = 17770
$ws.Range("M77").Value = -185197032
$ws.Range("N77").Value = -26506
# row 102
$ws.Range("H102").Value = 1368.25
$ws.Range("I102").Value = 1368.25
$ws.Range("K102").Value = 1368.25
$ws.Range("M102").Value = 253.75
# row 110
$ws.Range("H110").Value = 7901.3
$ws.Range("I110").Value = 6074.2
$ws.Range("J110").Value = 13382.6
$ws.Range("K110").Value = 6074.2
$ws.Range("L110").Value = 13382.6
$ws.Range("M110").Value = -4029.2
$ws.Range("N110").Value = -17472.6
# row 122
$ws.Range("H122").Value = 3619.4707
$ws.Range("I122").Value = 3022.08
$ws.Range("K122").Value = 9066.24
$ws.Range("M122").Value = -6616.24
# row 124
$ws.Range("H124").Value = 59962.5
$ws.Range("J124").Value = 59962.5
$ws.Range("L124").Value = 59962.5
$ws.Range("N124").Value = -69782.5
# row 132
$ws.Range("H132").Value = 3646.8
$ws.Range("I132").Value = 3272.0833
$ws.Range("J132").Value = 3992.6924
$ws.Range("K132").Value = 9816.249899999999
$ws.Range("L132").Value = 11978.0772
$ws.Range("M132").Value = -7286.249899999999
$ws.Range("N132").Value = -17038.0772
# row 136
$ws.Range("H136").Value = 6784.091
$ws.Range("I136").Value = 6292.2
$ws.Range("J136").Value = 7194
$ws.Range("K136").Value = 18876.6
$ws.Range("L136").Value = 21582
$ws.Range("M136").Value = -16326.6
$ws.Range("N136").Value = -26682

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 20
$ws.Range("H20").Value = 5470.478
$ws.Range("I20").Value = 4920.8887
$ws.Range("J20").Value = 7449
$ws.Range("K20").Value = 4920.8887
$ws.Range("L20").Value = 7449
$ws.Range("M20").Value = -4673.8887
$ws.Range("N20").Value = -7943
# row 134
$ws.Range("H134").Value = 1887.5416
$ws.Range("I134").Value = 1194.8636
$ws.Range("J134").Value = 9507
$ws.Range("K134").Value = 3584.5908
$ws.Range("L134").Value = 28521
$ws.Range("M134").Value = -1049.5908
$ws.Range("N134").Value = -33591

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 22
$ws.Range("H22").Value = 10162
$ws.Range("I22").Value = 9986.5
$ws.Range("K22").Value = 9986.5
$ws.Range("M22").Value = -9636.5
# row 38
$ws.Range("H38").Value = 1200
$ws.Range("J38").Value = 1200
$ws.Range("L38").Value = 1200
$ws.Range("N38").Value = -1954
# row 46
$ws.Range("H46").Value = 1200
$ws.Range("J46").Value = 1200
$ws.Range("L46").Value = 1200
$ws.Range("N46").Value = -1622
# row 132
$ws.Range("H132").Value = 4257.35
$ws.Range("I132").Value = 3473.8572
$ws.Range("J132").Value = 6085.5
$ws.Range("K132").Value = 10421.5716
$ws.Range("L132").Value = 18256.5
$ws.Range("M132").Value = -7891.571599999999
$ws.Range("N132").Value = -23316.5

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 93
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("M93").ClearContents()
# row 126
$ws.Range("H126").Value = 41669616
$ws.Range("I126").Value = 2233
$ws.Range("K126").Value = 6699
$ws.Range("M126").Value = -1759
# row 131
$ws.Range("H131").Value = 6265339.5
$ws.Range("J131").Value = 4421020
$ws.Range("L131").Value = 13263060
$ws.Range("N131").Value = -13273140

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 22
$ws.Range("H22").Value = 2281.8
$ws.Range("I22").Value = 2281.8
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 2281.8
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -1752.8
$ws.Range("N22").ClearContents()
# row 70
$ws.Range("H70").Value = 11908.454
$ws.Range("I70").Value = 9713.571
$ws.Range("K70").Value = 9713.571
$ws.Range("M70").Value = -9443.571
# row 73
$ws.Range("H73").Value = 11908.454
$ws.Range("I73").Value = 9713.571
$ws.Range("K73").Value = 9713.571
$ws.Range("M73").Value = -8777.571
# row 80
$ws.Range("H80").Value = 719641.1
$ws.Range("I80").Value = 836413.7
$ws.Range("J80").Value = 19006
$ws.Range("K80").Value = 836413.7
$ws.Range("L80").Value = 19006
$ws.Range("M80").Value = -835415.7
$ws.Range("N80").Value = -21002
# row 83
$ws.Range("H83").Value = 719641.1
$ws.Range("I83").Value = 836413.7
$ws.Range("J83").Value = 19006
$ws.Range("K83").Value = 4182068.5
$ws.Range("L83").Value = 95030
$ws.Range("M83").Value = -4177076.5
$ws.Range("N83").Value = -105014
# row 122
$ws.Range("H122").Value = 5599.88
$ws.Range("I122").Value = 3616.5557
$ws.Range("J122").Value = 7928.1304
$ws.Range("K122").Value = 10849.6671
$ws.Range("L122").Value = 23784.3912
$ws.Range("M122").Value = -8399.667099999999
$ws.Range("N122").Value = -28684.3912
# row 132
$ws.Range("H132").Value = 3657.8333
$ws.Range("I132").Value = 2769.0667
$ws.Range("K132").Value = 8307.2001
$ws.Range("M132").Value = -5777.2001

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 25
$ws.Range("H25").Value = 4000
$ws.Range("I25").Value = 4000
$ws.Range("K25").Value = 4000
$ws.Range("M25").Value = -3770
# row 32
$ws.Range("H32").Value = 3000
$ws.Range("I32").Value = 3000
$ws.Range("K32").Value = 3000
$ws.Range("M32").Value = -2683
# row 122
$ws.Range("H122").Value = 405360.7
$ws.Range("I122").Value = 804120.8
$ws.Range("K122").Value = 2412362.4
$ws.Range("M122").Value = -2409912.4
# row 132
$ws.Range("H132").Value = 10327.182
$ws.Range("I132").Value = 9780.4
$ws.Range("J132").Value = 10782.833
$ws.Range("K132").Value = 29341.2
$ws.Range("L132").Value = 32348.499
$ws.Range("M132").Value = -26811.2
$ws.Range("N132").Value = -37408.499
# row 136
$ws.Range("H136").Value = 2584.6309
$ws.Range("I136").Value = 1672.4706
$ws.Range("J136").Value = 5907.5
$ws.Range("K136").Value = 5017.4118
$ws.Range("L136").Value = 17722.5
$ws.Range("M136").Value = -2467.4118
$ws.Range("N136").Value = -22822.5
# row 140
$ws.Range("H140").Value = 77399.5
$ws.Range("J140").Value = 77399.5
$ws.Range("L140").Value = 77399.5
$ws.Range("N140").Value = -87759.5

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 122
$ws.Range("H122").Value = 3662.9412
$ws.Range("I122").Value = 2049.6365
$ws.Range("K122").Value = 6148.9095
$ws.Range("M122").Value = -3698.9095
# row 126
$ws.Range("H126").Value = 3839.3125
$ws.Range("I126").Value = 4819.727
$ws.Range("K126").Value = 14459.181
$ws.Range("M126").Value = -11989.181
# row 132
$ws.Range("H132").Value = 6103
$ws.Range("I132").Value = 2170
$ws.Range("K132").Value = 6510
$ws.Range("M132").Value = -3980
